$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.042.65'
$ws.Range('E2').Value = '  -0.62%  '
$ws.Range('D3').Value = '1.653.68'
$ws.Range('E3').Value = '  -0.04%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.32%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.10'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.63%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5217'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.35%  '
$ws.Range('E7').Value = '  -0.25%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2617'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.16%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06271'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.59%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.58'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.19%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07742'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.08%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.463'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.95%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.643.65'
$ws.Range('E13').Value = '  -0.67%  '
$ws.Range('D14').Value = '1.880.09'
$ws.Range('E14').Value = '  +0.04%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5426'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.50%  '
$ws.Range('D16').Value = '0.0₅8103'
$ws.Range('E16').Value = '  -1.11%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.97'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.38%  '
$ws.Range('D18').Value = '26.040.55'
$ws.Range('E18').Value = '  -0.60%  '
$ws.Range('E19').Value = '  -0.30%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.574'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.20%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '191.67'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.37%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.02'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.39%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.977'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '138.36'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.16%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1236'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.41%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.258'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.33%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.14'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.50%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.405'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.77%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05962'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.62%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.273'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.82%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.507'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.13%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.241'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.48%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.564'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.29%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9495'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.52%  '
$ws.Range('E36').Value = '  -0.11%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.750'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.71%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5689'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.98%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01599'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.21%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.887'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.09%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8442'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.13%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.002'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.21%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '100.89'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.16%  '
$ws.Range('D44').Value = '1.007.29'
$ws.Range('E44').Value = '  -4.76%  '
$ws.Range('D45').Value = '1.794.42'
$ws.Range('E45').Value = '  -0.01%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '56.71'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.01%  '
$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').Value = '0.0₈107'
$ws.Range('E47').Value = '  +1.20%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.001'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.34%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.965'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.94%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4301'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.64%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.475'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.57%  '
